$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the whole "Meta description" paragraph that currently sits
#    right after the title (Heading1) paragraph.
# ------------------------------------------------------------------
$metaRange = $d.Content.Duplicate
$foundMeta = $metaRange.Find.Execute("Meta description: Read our review of Cyberslot Megaclusters, the future-themed slot machine with mini-grids and a wandering Wild symbol. Play it free and experience unique gameplay mechanics.", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundMeta) {
    $metaPara = $metaRange.Paragraphs.Item(1)
    # Extend the delete range by one character so the trailing paragraph
    # mark is removed too (otherwise an empty paragraph is left behind).
    $metaParaFull = $d.Range($metaPara.Range.Start, $metaPara.Range.End + 1)
    $metaParaFull.Delete()
}

# ------------------------------------------------------------------
# 2. Insert a new bold paragraph ("Play Cyberslot Megaclusters Free:
#    Unique Gameplay Mechanics") right before the final "Prompt: ..."
#    paragraph.
# ------------------------------------------------------------------
$anchorRange = $d.Content.Duplicate
$foundAnchor = $anchorRange.Find.Execute("High volatility that might not be suitable for casual players", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundAnchor) {
    $insertionPoint = $d.Range($anchorRange.End - 1, $anchorRange.End - 1)
    $newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Cyberslot Megaclusters Free: Unique Gameplay Mechanics</w:t></w:r></w:p>'
    [void]$insertionPoint.InsertXML($newParaXml)
}

# ------------------------------------------------------------------
# 3. Replace the text of the final paragraph (formerly the "Prompt:"
#    image-generation prompt) with the meta description copy, keeping
#    its existing italic run formatting.
# ------------------------------------------------------------------
$promptRange = $d.Content.Duplicate
$foundPrompt = $promptRange.Find.Execute("Prompt: Create a colorful and eye-catching feature image in a cartoon style for Cyberslot Megaclusters. The image should prominently feature a happy Maya warrior wearing glasses, in line with the futuristic and technology-themed game. The image should convey the game's mini-grid mechanic and use a bright color palette to reflect the game's simple yet modern aesthetic. Please include the game title and any additional elements that you feel would enhance the image's appeal and accurately represent the game's features.", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundPrompt) {
    $promptRange.Text = "Read our review of Cyberslot Megaclusters, the future-themed slot machine with mini-grids and a wandering Wild symbol. Play it free and experience unique gameplay mechanics."
}

Write-Host "Edit complete"
